$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark a few items as "Done = Y" (with the built-in "Good" cell style) ---
$doneCells = @("H6", "H12", "H13", "H15")
foreach ($cellRef in $doneCells) {
    $ws.Range($cellRef).Value = "Y"
    $ws.Range($cellRef).Style = "Good"
}

# --- Move the "Shipment fee Ebay" amount from the P U column to the Price column ---
$ws.Range("F24").Value = $ws.Range("E24").Value()
$ws.Range("E24").ClearContents()

# --- Add a new line for the openEVSE shipment fee ---
$ws.Range("D25").Value = "Shipment fee openEVSE"
$ws.Range("F25").Formula = "=26.37/1.2"

# --- Add a grand total row ---
$ws.Range("F27").Formula = "=F23+F24+F25"

# --- Update the selection to match the latest edit ---
[void]$ws.Range("G16:G18").Select()
